# Update the "Sparsity coeff" (column L) values to reflect the new
# definition: sparsity coefficient is now the AVERAGE sparsity of the
# coefficient matrix, with adaptive sparsity introduced for the @dch
# and @dchperceptron algorithm rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L4").Value  = 0.70923000000000003
$ws.Range("L5").Value  = 0.70611999999999997
$ws.Range("L13").Value = 0.18210000000000001
$ws.Range("L14").Value = 0.18221000000000001
$ws.Range("L22").Value = 0.32434000000000002
$ws.Range("L23").Value = 0.32425999999999999
$ws.Range("L30").Value = 5.4715800000000003
$ws.Range("L31").Value = 5.4707499999999998
$ws.Range("L37").Value = 5.0435400000000001
$ws.Range("L38").Value = 5.0436800000000002

# Reflect the author's last on-screen position: scrolled down/selected
# a different cell than before.
$ws.Select()
try {
    $excel.ActiveWindow.ScrollRow = 16
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}
$ws.Range("I43").Select()
